$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.537.09"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -2.46%  "
$c.Style = "Normal"

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.582.54"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -3.08%  "
$c.Style = "Normal"

# Row 4
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c.Style = "Normal"

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "210.70"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -2.76%  "
$c.Style = "Normal"

# Row 6
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -1.92%  "
$c.Style = "Normal"

# Row 7
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +0.08%  "
$c.Style = "Normal"

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.249"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -2.36%  "
$c.Style = "Normal"

# Row 9
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -1.10%  "
$c.Style = "Normal"

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "19.48"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -4.04%  "
$c.Style = "Normal"

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0832"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -1.98%  "
$c.Style = "Normal"

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.802.64"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -3.17%  "
$c.Style = "Normal"

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.580.14"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -3.12%  "
$c.Style = "Normal"

# Row 14
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -1.78%  "
$c.Style = "Normal"

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.528"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -2.99%  "
$c.Style = "Normal"

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "64.12"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -1.08%  "
$c.Style = "Normal"

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "26.554.90"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -2.34%  "
$c.Style = "Normal"

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.0₃0728"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -0.78%  "
$c.Style = "Normal"

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "208.51"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -3.42%  "
$c.Style = "Normal"

# Row 20
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +0.28%  "
$c.Style = "Normal"

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.69"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -3.09%  "
$c.Style = "Normal"

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.25"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -3.55%  "
$c.Style = "Normal"

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.40"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -4.02%  "
$c.Style = "Normal"

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "8.89"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -2.10%  "
$c.Style = "Normal"

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "146.47"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -1.34%  "
$c.Style = "Normal"

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.43"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +1.92%  "
$c.Style = "Normal"

# Row 27
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -0.11%  "
$c.Style = "Normal"

# Row 28
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -4.62%  "
$c.Style = "Normal"

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "15.27"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -1.98%  "
$c.Style = "Normal"

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0501"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -0.70%  "
$c.Style = "Normal"

# Row 31
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -2.38%  "
$c.Style = "Normal"

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.26"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -3.79%  "
$c.Style = "Normal"

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.660"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +22.47%  "
$c.Style = "Normal"

# Row 34
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -2.69%  "
$c.Style = "Normal"

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.309.72"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -0.59%  "
$c.Style = "Normal"

# Row 36
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.50"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -3.78%  "
$c.Style = "Normal"

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.43"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -0.65%  "
$c.Style = "Normal"

# Row 38
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -1.22%  "
$c.Style = "Normal"

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.821"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -3.34%  "
$c.Style = "Normal"

# Row 40
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +0.17%  "
$c.Style = "Normal"

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.784"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -2.56%  "
$c.Style = "Normal"

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.29"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +0.95%  "
$c.Style = "Normal"

# Row 43
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.17"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -4.36%  "
$c.Style = "Normal"

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "62.86"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -1.46%  "
$c.Style = "Normal"

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.716.76"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -2.99%  "
$c.Style = "Normal"

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "88.88"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  -2.23%  "
$c.Style = "Normal"

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.62"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +0.67%  "
$c.Style = "Normal"

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.827"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +3.68%  "
$c.Style = "Normal"

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0984"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +3.62%  "
$c.Style = "Normal"

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0506"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -1.97%  "
$c.Style = "Normal"

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.45"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -1.44%  "
$c.Style = "Normal"
